# Auto-generated Excel COM-interop script to apply Marilith_Profits market data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 277.2857
$ws.Range("I6").Value = 80.333336
$ws.Range("J6").Value = 425
$ws.Range("K6").Value = 241.000008
$ws.Range("L6").Value = 1275
$ws.Range("M6").Value = -129.000008
$ws.Range("N6").Value = -1499
# Row 33
$ws.Range("H33").Value = 649.25
$ws.Range("I33").Value = 287.8889
$ws.Range("K33").Value = 287.8889
$ws.Range("M33").Value = -58.88889999999998
# Row 38
$ws.Range("H38").Value = 28.875
$ws.Range("I38").Value = 28.875
$ws.Range("K38").Value = 86.625
$ws.Range("M38").Value = 285.375
# Row 39
$ws.Range("H39").Value = 225.53334
$ws.Range("I39").Value = 170.21428
$ws.Range("K39").Value = 510.64284
$ws.Range("M39").Value = -214.64284
# Row 51
$ws.Range("H51").Value = 4997.5
$ws.Range("I51").Value = 4998
$ws.Range("K51").Value = 4998
$ws.Range("M51").Value = -4514
# Row 70
$ws.Range("H70").Value = 1068972
$ws.Range("I70").Value = 3375981.2
$ws.Range("J70").Value = 4198.5386
$ws.Range("K70").Value = 10127943.6
$ws.Range("L70").Value = 12595.6158
$ws.Range("M70").Value = -10127673.6
$ws.Range("N70").Value = -13135.6158
# Row 73
$ws.Range("H73").Value = 1068972
$ws.Range("I73").Value = 3375981.2
$ws.Range("J73").Value = 4198.5386
$ws.Range("K73").Value = 10127943.6
$ws.Range("L73").Value = 12595.6158
$ws.Range("M73").Value = -10127007.6
$ws.Range("N73").Value = -14467.6158
# Row 99
$ws.Range("H99").Value = 493.6
$ws.Range("I99").Value = 493.6
$ws.Range("K99").Value = 1480.8
$ws.Range("M99").Value = 17.19999999999982
# Row 137
$ws.Range("H137").Value = 2713.7334
$ws.Range("I137").Value = 1897.8889
$ws.Range("J137").Value = 3937.5
$ws.Range("K137").Value = 5693.6667
$ws.Range("L137").Value = 11812.5
$ws.Range("M137").Value = -3143.6667
$ws.Range("N137").Value = -16912.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1998
$ws.Range("I2").Value = 1998
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1998
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1885
# Row 32
$ws.Range("H32").Value = 9219.762000000001
$ws.Range("I32").Value = 9219.762000000001
$ws.Range("K32").Value = 9219.762000000001
$ws.Range("M32").Value = -8932.762000000001
# Row 61
$ws.Range("H61").Value = 2331.3333
$ws.Range("I61").Value = 1997.5
$ws.Range("K61").Value = 1997.5
$ws.Range("M61").Value = -1785.5
# Row 74
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2126
# Row 77
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10632
# Row 116
$ws.Range("H116").Value = 1998
$ws.Range("I116").Value = 1998
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1998
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 296
# Row 132
$ws.Range("H132").Value = 4016.9583
$ws.Range("I132").Value = 3713.348
$ws.Range("K132").Value = 11140.044
$ws.Range("M132").Value = -8610.044
# Row 136
$ws.Range("H136").Value = 2331.3333
$ws.Range("I136").Value = 1997.5
$ws.Range("K136").Value = 5992.5
$ws.Range("M136").Value = -3442.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1998
$ws.Range("I3").Value = 1998
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1998
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1884
# Row 20
$ws.Range("H20").Value = 3786.5
$ws.Range("J20").Value = 5473
$ws.Range("L20").Value = 5473
$ws.Range("N20").Value = -5967
# Row 107
$ws.Range("H107").Value = 1478.625
$ws.Range("I107").Value = 956
$ws.Range("K107").Value = 956
$ws.Range("M107").Value = 964

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 426.46155
$ws.Range("I5").Value = 231
$ws.Range("J5").Value = 739.2
$ws.Range("K5").Value = 231
$ws.Range("L5").Value = 739.2
$ws.Range("M5").Value = -119
$ws.Range("N5").Value = -963.2
# Row 31
$ws.Range("H31").Value = 5451.8
$ws.Range("I31").Value = 5999
$ws.Range("J31").Value = 5315
$ws.Range("K31").Value = 5999
$ws.Range("L31").Value = 5315
$ws.Range("M31").Value = -5704
$ws.Range("N31").Value = -5905
# Row 34
$ws.Range("H34").Value = 5451.8
$ws.Range("I34").Value = 5999
$ws.Range("J34").Value = 5315
$ws.Range("K34").Value = 5999
$ws.Range("L34").Value = 5315
$ws.Range("M34").Value = -5797
$ws.Range("N34").Value = -5719
# Row 41
$ws.Range("H41").Value = 14284.667
$ws.Range("J41").Value = 21250
$ws.Range("L41").Value = 21250
$ws.Range("N41").Value = -22106
# Row 59
$ws.Range("H59").Value = 27991.545
$ws.Range("I59").Value = 22151.166
$ws.Range("K59").Value = 22151.166
$ws.Range("M59").Value = -21006.166
# Row 68
$ws.Range("H68").Value = 51426.8
$ws.Range("I68").Value = 34268
$ws.Range("J68").Value = 53333.332
$ws.Range("K68").Value = 34268
$ws.Range("L68").Value = 53333.332
$ws.Range("M68").Value = -33519
$ws.Range("N68").Value = -54831.332
# Row 71
$ws.Range("H71").Value = 51426.8
$ws.Range("I71").Value = 34268
$ws.Range("J71").Value = 53333.332
$ws.Range("K71").Value = 102804
$ws.Range("L71").Value = 159999.996
$ws.Range("M71").Value = -99060
$ws.Range("N71").Value = -167487.996

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 999.75
$ws.Range("I132").Value = 999.75
$ws.Range("K132").Value = 2999.25
$ws.Range("M132").Value = -469.25

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 799
$ws.Range("I13").Value = 798
$ws.Range("K13").Value = 798
$ws.Range("M13").Value = -658
# Row 16
$ws.Range("H16").Value = 3841.625
$ws.Range("I16").Value = 4176.143
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 4176.143
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -4006.143
$ws.Range("N16").Value = -1840
# Row 136
$ws.Range("H136").Value = 3251
$ws.Range("I136").Value = 3251
$ws.Range("K136").Value = 9753
$ws.Range("M136").Value = -7203

$ws = $wb.Worksheets.Item("WVR")
# Row 117
$ws.Range("H117").Value = 39000
$ws.Range("J117").Value = 39000
$ws.Range("L117").Value = 39000
$ws.Range("N117").Value = -48178
